$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SV_calls")
Write-Host $ws.Name
